$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 27.12.2021 change in fun createNewType, createVersion, CreateRevision
# The BOM generator now produces a new ParentId ("B12UAZ1111212") for the
# parent assembly described on this sheet. Update every child row's
# ParentId (column A) to the newly generated value.
$newParentId = "B12UAZ1111212"

$ws.Range("A2").Value = $newParentId
$ws.Range("A3").Value = $newParentId
$ws.Range("A4").Value = $newParentId
$ws.Range("A5").Value = $newParentId
$ws.Range("A6").Value = $newParentId
